$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11, shifting the existing record set (rows 11-108)
# down by one row (they become rows 12-109), and fill the new row 11 with a
# new weekly price record for Poroto granado.
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = 2
$ws.Range("B11").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C11").Value = 'Coquimbo'
$ws.Range("D11").Value = 44959
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 100112030
$ws.Range("G11").Value = 'Poroto granado'
$ws.Range("H11").Value = 'Sin especificar'
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 600
$ws.Range("K11").Value = 23000
$ws.Range("L11").Value = 25000
$ws.Range("M11").Value = 24000
$ws.Range("N11").Value = '$/malla 25 kilos'
$ws.Range("O11").Value = 'Provincia de Limarí'
$ws.Range("P11").Value = 960
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = 'Hortaliza'
